# Update "Overview" income-statement sheet:
#  - drop the oldest reporting period (12 months ended 1396/12) and shift
#    every period one column to the left (D<-E, E<-F, F<-G, G<-H)
#  - append a brand-new period column H for the 12 months ended 1401/12,
#    together with its figures
#  - refresh the "تاریخ انتشار" (publish date) row: the amended/re-issued
#    date for the 1400/12 statement (now column G) and the first publish
#    date for the new 1401/12 statement (column H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Row, $Values) {
    $cols = @("D", "E", "F", "G", "H")
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $addr = "$($cols[$i])$Row"
        $ws.Range($addr).Value = $Values[$i]
    }
}

# Row 8: period headers
Set-RowValues 8 @(
    "12 ماهه منتهی به 1397/12",
    "12 ماهه منتهی به 1398/12",
    "12 ماهه منتهی به 1399/12",
    "12 ماهه منتهی به 1400/12",
    "12 ماهه منتهی به 1401/12"
)

# Row 9: publish dates
Set-RowValues 9 @(
    "1399-02-31 (10)",
    "1400-02-30 (8)",
    "1401-02-31 (8)",
    "1402-02-10 (7)",
    "1402-02-10"
)

# Row 11: فروش (Sales)
Set-RowValues 11 @(94001, 110105, 97997, 199860, 224892)

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
Set-RowValues 12 @(-36503, -43217, -32596, -71505, -70152)

# Row 13: سود (زیان) ناخالص (Gross profit)
Set-RowValues 13 @(57498, 66887, 65400, 128355, 154740)

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
Set-RowValues 14 @(-8962, -6323, -5365, -7916, -14043)

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net)
Set-RowValues 16 @(3511, 2608, 743, 1793, 5079)

# Row 17: سود (زیان) عملیاتی (Operating profit)
Set-RowValues 17 @(52047, 63173, 60779, 122232, 145776)

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net)
Set-RowValues 19 @(6736, 8365, 6502, 9214, 11437)

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing operations)
Set-RowValues 20 @(58783, 71538, 67281, 131446, 157213)

# Row 21: مالیات (Tax)
Set-RowValues 21 @(-3528, -3328, -3739, -8442, -17026)

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing operations)
Set-RowValues 22 @(55255, 68210, 63542, 123004, 140187)

# Row 24: سود (زیان) خالص (Net profit)
Set-RowValues 24 @(55255, 68210, 63542, 123004, 140187)

# Row 26: سود هر سهم پس از کسر مالیات (EPS after tax)
Set-RowValues 26 @(17694, 13953, 7916, 6784, 5072)
